$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header N1: rename model label "Modelo 3 - log" -> "Modelo 3 - Forest"
$ws.Range("N1").Value = "Modelo 3 - Forest"

# Updated inferred quantities (N column) - new random seed / criteria results
$ws.Range("N4").Value = 46
$ws.Range("N6").Value = 40
$ws.Range("N9").Value = 30
$ws.Range("N11").Value = 84
$ws.Range("N13").Value = 48
$ws.Range("N15").Value = 48
$ws.Range("N17").Value = 212
$ws.Range("N18").Value = 46
$ws.Range("N19").Value = 40

# Rows 16 and 17 no longer have comparison formulas (P,Q,R cleared but keep formatting)
$ws.Range("P16").ClearContents()
$ws.Range("Q16").ClearContents()
$ws.Range("R16").ClearContents()

$ws.Range("P17").ClearContents()
$ws.Range("Q17").ClearContents()
$ws.Range("R17").ClearContents()

# Selection moved
$ws.Range("U18").Select()
